$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "51.983.08"
Set-TextCell $ws.Range("E2") "  -0.23%  "
Set-TextCell $ws.Range("D3") "2.792.82"
Set-TextCell $ws.Range("E3") "  -1.77%  "
Set-TextCell $ws.Range("E4") "  -0.03%  "
Set-TextCell $ws.Range("D5") "358.83"
Set-TextCell $ws.Range("E5") "  -1.04%  "
Set-TextCell $ws.Range("D6") "109.84"
Set-TextCell $ws.Range("E6") "  -3.38%  "
Set-TextCell $ws.Range("D7") "0.558"
Set-TextCell $ws.Range("E7") "  -1.59%  "
Set-TextCell $ws.Range("D8") "1.00"
Set-TextCell $ws.Range("E8") "  +0.01%  "
Set-TextCell $ws.Range("E9") "  -2.24%  "
Set-TextCell $ws.Range("D10") "39.80"
Set-TextCell $ws.Range("E10") "  -4.69%  "
Set-TextCell $ws.Range("D11") "0.0847"
Set-TextCell $ws.Range("E11") "  -1.69%  "
Set-TextCell $ws.Range("E12") "  +1.19%  "
Set-TextCell $ws.Range("D13") "19.42"
Set-TextCell $ws.Range("E13") "  -2.96%  "
Set-TextCell $ws.Range("D14") "7.55"
Set-TextCell $ws.Range("E14") "  -3.05%  "
Set-TextCell $ws.Range("D15") "3.228.67"
Set-TextCell $ws.Range("E15") "  -1.96%  "
Set-TextCell $ws.Range("D16") "2.803.92"
Set-TextCell $ws.Range("E16") "  -1.42%  "
Set-TextCell $ws.Range("D17") "0.937"
Set-TextCell $ws.Range("E17") "  +3.49%  "
Set-TextCell $ws.Range("D18") "51.940.82"
Set-TextCell $ws.Range("E18") "  -0.07%  "
Set-TextCell $ws.Range("D19") "7.46"
Set-TextCell $ws.Range("E19") "  +0.12%  "
Set-TextCell $ws.Range("E20") "  -2.64%  "
Set-TextCell $ws.Range("D21") "13.06"
Set-TextCell $ws.Range("E21") "  -3.65%  "
Set-TextCell $ws.Range("E22") "  -1.65%  "
Set-TextCell $ws.Range("D23") "70.18"
Set-TextCell $ws.Range("E23") "  +0.05%  "
Set-TextCell $ws.Range("D24") "270.21"
Set-TextCell $ws.Range("E24") "  +1.16%  "
Set-TextCell $ws.Range("D25") "2.75"
Set-TextCell $ws.Range("E25") "  -3.73%  "
Set-TextCell $ws.Range("D26") "26.53"
Set-TextCell $ws.Range("E26") "  -2.14%  "
Set-TextCell $ws.Range("E27") "  -0.04%  "
Set-TextCell $ws.Range("D28") "0.165"
Set-TextCell $ws.Range("E28") "  +18.26%  "
Set-TextCell $ws.Range("D29") "10.24"
Set-TextCell $ws.Range("E29") "  -1.84%  "
Set-TextCell $ws.Range("D30") "2.21"
Set-TextCell $ws.Range("E30") "  -1.65%  "
Set-TextCell $ws.Range("D31") "0.0472"
Set-TextCell $ws.Range("E31") "  +5.55%  "
Set-TextCell $ws.Range("D32") "52.07"
Set-TextCell $ws.Range("E32") "  -2.84%  "
Set-TextCell $ws.Range("D33") "33.70"
Set-TextCell $ws.Range("E33") "  -1.11%  "
Set-TextCell $ws.Range("D34") "5.73"
Set-TextCell $ws.Range("E34") "  -2.85%  "
Set-TextCell $ws.Range("D35") "0.0840"
Set-TextCell $ws.Range("E35") "  +0.18%  "
Set-TextCell $ws.Range("D36") "5.20"
Set-TextCell $ws.Range("E36") "  -1.97%  "
Set-TextCell $ws.Range("E37") "  -0.12%  "
Set-TextCell $ws.Range("D38") "18.79"
Set-TextCell $ws.Range("E38") "  +2.52%  "
Set-TextCell $ws.Range("D39") "3.21"
Set-TextCell $ws.Range("E39") "  -2.86%  "
Set-TextCell $ws.Range("D40") "2.00"
Set-TextCell $ws.Range("E40") "  -4.14%  "
Set-TextCell $ws.Range("D41") "2.56"
Set-TextCell $ws.Range("E41") "  -0.12%  "
Set-TextCell $ws.Range("E42") "  -1.52%  "
Set-TextCell $ws.Range("D43") "2.25"
Set-TextCell $ws.Range("E43") "  -0.55%  "
Set-TextCell $ws.Range("D44") "119.63"
Set-TextCell $ws.Range("E44") "  -6.36%  "
Set-TextCell $ws.Range("D45") "21.81"
Set-TextCell $ws.Range("E45") "  -10.31%  "
Set-TextCell $ws.Range("D46") "2.081.17"
Set-TextCell $ws.Range("E46") "  -1.87%  "
Set-TextCell $ws.Range("D47") "3.24"
Set-TextCell $ws.Range("E47") "  -4.36%  "
Set-TextCell $ws.Range("D48") "2.22"
Set-TextCell $ws.Range("E48") "  -1.73%  "
Set-TextCell $ws.Range("D49") "5.83"
Set-TextCell $ws.Range("E49") "  +0.01%  "
Set-TextCell $ws.Range("D50") "0.955"
Set-TextCell $ws.Range("E50") "  -5.11%  "
Set-TextCell $ws.Range("D51") "8.88"
Set-TextCell $ws.Range("E51") "  -1.60%  "
